$d = $word.ActiveDocument

# objetivos_pt
$old = "Levar aos estudantes conhecimentos básicos sobre:a) Biotecnologia: enfatizando o emprego dos processos bioquímicos relevantes para as diferentes áreas de aplicação da biotecnologia.b) Processos Fermentativos / Enzimáticos: compreendendo conceitos, características e etapas de desenvolvimento.c) Bioquímica das fermentações: focando nas principais rotas metabólicas utilizadas por microrganismos de interesse industrial"
$new = "Levar aos estudantes conhecimentos básicos sobre:^la) Biotecnologia: enfatizando o emprego dos processos bioquímicos relevantes para as diferentes áreas de aplicação da biotecnologia.^lb) Processos Fermentativos / Enzimáticos: compreendendo conceitos, características e etapas de desenvolvimento.^lc) Bioquímica das fermentações: focando nas principais rotas metabólicas utilizadas por microrganismos de interesse industrial"
$found = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
if (-not $found) { Write-Output "NOT FOUND: objetivos_pt" }

# objetivos_en
$old = "The aims of this course are focused on some relevant issues regarding biotechnology(field of applications); fermentative and enzymatic processes; biochemistry of thefermentations (metabolic pathways of industrial interest); fermentative processes ofindustrial interest"
$new = "The aims of this course are focused on some relevant issues regarding biotechnology^l(field of applications); fermentative and enzymatic processes; biochemistry of the^lfermentations (metabolic pathways of industrial interest); fermentative processes of^lindustrial interest"
$found = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
if (-not $found) { Write-Output "NOT FOUND: objetivos_en" }

# programa_en
$old = "1. Biotechnology: concepts, application areas, multidisciplinary characteristic andexamples of biotechnological products and processes.2. Fermentative processes: concept, enzymatic and fermentative processes, steps offermentative process (downstream x upstream). Fermentative process modes: a) batchand fed-batch fermentation, semi continuous and, continuous processes; b) induced andspontaneous fermentation; c) semi solid fermentation; d) oxygen supply; e) submergedand in surface processes; f) kinetics of the product formation in relation to the primarymetabolism according to Gaden.3. Biochemistry of the fermentation: Fermentation – concepts, objectives, aerobic xanaerobic metabolisms; energy balance; preliminary steps of fermentation (extracellularhydrolysis and membrane permeability); metabolic pathways of industrial interest: a)EMP pathway; reactions and allosteric control; alcoholic fermentation, homolacticfermentation, acetone/butanol, mixed-acid and 2,3 butanediol; b) Fosfo-Ketolasepathway; heterolactic fermentation and c) Entner Doudoroff pathway: alcoholicfermentation by Zymmonas mobilis. Evaluation parameters of a fermentative process:yield, fermentation efficiency and productivity. Highlights of some processes ofindustrial interest, such as cocoa processing, ethanol production, fermented food andothers."
$new = "1. Biotechnology: concepts, application areas, multidisciplinary characteristic and^lexamples of biotechnological products and processes.^l^l2. Fermentative processes: concept, enzymatic and fermentative processes, steps of^lfermentative process (downstream x upstream). Fermentative process modes: a) batch^land fed-batch fermentation, semi continuous and, continuous processes; b) induced and^lspontaneous fermentation; c) semi solid fermentation; d) oxygen supply; e) submerged^land in surface processes; f) kinetics of the product formation in relation to the primary^lmetabolism according to Gaden.^l^l3. Biochemistry of the fermentation: Fermentation – concepts, objectives, aerobic x^lanaerobic metabolisms; energy balance; preliminary steps of fermentation (extracellular^lhydrolysis and membrane permeability); metabolic pathways of industrial interest: a)^lEMP pathway; reactions and allosteric control; alcoholic fermentation, homolactic^lfermentation, acetone/butanol, mixed-acid and 2,3 butanediol; b) Fosfo-Ketolase^lpathway; heterolactic fermentation and c) Entner Doudoroff pathway: alcoholic^lfermentation by Zymmonas mobilis. Evaluation parameters of a fermentative process:^lyield, fermentation efficiency and productivity. Highlights of some processes of^lindustrial interest, such as cocoa processing, ethanol production, fermented food and^lothers."
$found = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
if (-not $found) { Write-Output "NOT FOUND: programa_en" }

# bibliografia
$old = "1. AMERINE, M.A, OUGH,C.S., Methods for analysis of musts and wines. New York: John Wiley & Sons, 1980. 2. AMORIM, H.V., Fermentação Alcoólica ciência e tecnologia. Piracicaba: Fermentec,2006.3. BORZANI, W., SCHMIDELL, W., LIMA, U.A., AQUARONE, E. Série de Biotecnologia Vol. 1 – Fundamentos e Vol. 4 Processos Fermentativos e Enzimáticos. São Paulo: Ed.Edgard Blucher, 2020.4. EL-MANSI, E.M.T., BRYCE, C.E.A., DEMAIN, A.L., ALLMAN,A.R. Fermentation Microbiology and Biotechnology. 2ª Ed. New York: CRC Taylor & Francis, 2007.5. SILVA, N.; TANIWAKI, M H., SA, P. B. Z. R.  Fermentação e processos fermentativos – São Paulo: Tiki Books: The Good Food Institute Brasil, 2022. (Série Tecnológica das Proteínas Alternativas) E-Book: PDF, 40 p.; IL6. BASTOS, R. G.; Tecnologia das fermentações: fundamentos de Bioprocessos. -- São Carlos :  EdUFSCar, 2010. 162 p. -- (Coleção UAB-UFSCar)."
$new = "1. AMERINE, M.A, OUGH,C.S., Methods for analysis of musts and wines. New York: John Wiley & Sons, 1980.^l ^l2. AMORIM, H.V., Fermentação Alcoólica ciência e tecnologia. Piracicaba: Fermentec,2006.^l^l3. BORZANI, W., SCHMIDELL, W., LIMA, U.A., AQUARONE, E. Série de Biotecnologia Vol. 1 – Fundamentos e Vol. 4 Processos Fermentativos e Enzimáticos. São Paulo: Ed.Edgard Blucher, 2020.^l^l4. EL-MANSI, E.M.T., BRYCE, C.E.A., DEMAIN, A.L., ALLMAN,A.R. Fermentation Microbiology and Biotechnology. 2ª Ed. New York: CRC Taylor & Francis, 2007.^l^l5. SILVA, N.; TANIWAKI, M H., SA, P. B. Z. R.  Fermentação e processos fermentativos – São Paulo: Tiki Books: The Good Food Institute Brasil, 2022. (Série Tecnológica das Proteínas Alternativas) E-Book: PDF, 40 p.; IL^l^l6. BASTOS, R. G.; Tecnologia das fermentações: fundamentos de Bioprocessos. -- São Carlos :  EdUFSCar, 2010. 162 p. -- (Coleção UAB-UFSCar)."
$found = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
if (-not $found) { Write-Output "NOT FOUND: bibliografia" }
